$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Size Behavior test block (rows 10-16) ---
$ws.Range("F10").Value2 = "aBagOfIntegers contains 10 integers"
$ws.Range("F11").Value2 = "aBagOfIntegers contains 10 integers"
$ws.Range("A12").Value2 = "Test ID: Size Behavior"
$ws.Range("F13").Value2 = "aBagOfIntegers contains 10 integers"
$ws.Range("B14").Value2 = "Input Data:"
$ws.Range("F14").Value2 = "aBagOfIntegers.size();"
$ws.Range("F15").Value2 = "`"There are now 10 items in your bag.`""
$ws.Range("F16").Value2 = "none"

# --- Count Behavior w/ duplicate test block (rows 17-21) ---
$ws.Range("A17").Value2 = "Test ID: Count Behavior w/ duplicate"
$ws.Range("F18").Value2 = "added another 9 to the bag"
$ws.Range("B19").Value2 = "Input Data:"
$ws.Range("F19").Value2 = "aBagOfIntegers.count(9);"
$ws.Range("F20").Value2 = "`"The number '9' is in this list 2 times.`""
$ws.Range("F21").Value2 = "none"

# --- Count Behavior w/o duplicate test block (rows 22-26) ---
$ws.Range("A22").Value2 = "Test ID: Count Behavior w/o duplicate"
$ws.Range("F23").Value2 = "there is only one "
$ws.Range("F24").Value2 = "aBagOfIntegers.count(1);"
$ws.Range("F25").Value2 = "`"The number '1' is in this list 1 times.`""

# --- Count Behavior w/ no instance test block (rows 27-31) ---
$ws.Range("A27").Value2 = "Test ID: Count Behavior w/ no instance"
$ws.Range("F28").Value2 = "there is none"
$ws.Range("F29").Value2 = "aBagOfIntegers.count(25);"
$ws.Range("B30").Value2 = "Expected Result"
$ws.Range("F30").Value2 = "`"The number '25' is in this list 0 times.`""

# Update selected/active cell to F30
$ws.Range("F30").Select()
